$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 332, pushing the existing row 332 (and all
# rows below it, through the old row 372) down by one, to 333..373.
$ws.Rows.Item(332).Insert()

# Populate the newly-inserted row 332 with the new weekly price record.
$ws.Cells.Item(332, 1).Value = 9
$ws.Cells.Item(332, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(332, 3).Value = "Metropolitana"
$ws.Cells.Item(332, 4).Value = 44946
$ws.Cells.Item(332, 5).Value = 13
$ws.Cells.Item(332, 6).Value = 300000001
$ws.Cells.Item(332, 7).Value = "Rabanito"
$ws.Cells.Item(332, 8).Value = "Sin especificar"
$ws.Cells.Item(332, 9).Value = "Primera"
$ws.Cells.Item(332, 10).Value = 7000
$ws.Cells.Item(332, 11).Value = 3000
$ws.Cells.Item(332, 12).Value = 3000
$ws.Cells.Item(332, 13).Value = 3000
$ws.Cells.Item(332, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(332, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(332, 16).Value = 30
$ws.Cells.Item(332, 17).Value = 100
$ws.Cells.Item(332, 18).Value = "Hortaliza"
